$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price (D) and Volume(1h) (E) columns store plain text that looks
# numeric (e.g. "244.12", "0.0002000", "-0.56%"), preserved exactly as
# typed (including trailing zeros). Force the cells to a text number
# format before writing so Excel does not silently convert them to
# numbers and strip formatting, then restore the default "Normal"
# style afterwards so no stray cell styles are introduced.
$textRange = $ws.Range("D2:E50")
$textRange.NumberFormat = "@"

# Row 2 - BNB
$ws.Range("D2").Value = "244.12"
$ws.Range("E2").Value = "-0.56%"

# Row 3 - OKB
$ws.Range("D3").Value = "26.58"
$ws.Range("E3").Value = "4.17%"

# Row 4 - HuobiToken
$ws.Range("D4").Value = "5.135"
$ws.Range("E4").Value = "0.00%"

# Row 5 - Cronos
$ws.Range("D5").Value = "0.05610"
$ws.Range("E5").Value = "0.48%"

# Row 6 - KuCoinToken
$ws.Range("D6").Value = "6.465"
$ws.Range("E6").Value = "-0.29%"

# Row 7 - MXToken
$ws.Range("E7").Value = "0.11%"

# Row 8 - FTXToken
$ws.Range("E8").Value = "-1.12%"

# Row 9 - WazirX
$ws.Range("D9").Value = "0.1329"
$ws.Range("E9").Value = "-0.48%"

# Row 10 - MandalaExchangeToken
$ws.Range("D10").Value = "0.06926"
$ws.Range("E10").Value = "-0.50%"

# Row 11 - BitrueCoin
$ws.Range("D11").Value = "0.02895"
$ws.Range("E11").Value = "1.23%"

# Row 12 - BitMartToken
$ws.Range("D12").Value = "0.09385"
$ws.Range("E12").Value = "0.07%"

# Row 13 - BitForexToken
$ws.Range("D13").Value = "0.001520"
$ws.Range("E13").Value = "0.33%"

# Row 14 - One
$ws.Range("D14").Value = "0.0005979"
$ws.Range("E14").Value = "0.52%"

# Row 15 - TigerCash
$ws.Range("D15").Value = "0.006152"
$ws.Range("E15").Value = "0.09%"

# Row 16 - LEO
$ws.Range("D16").Value = "3.648"
$ws.Range("E16").Value = "3.26%"

# Row 17 - GateToken
$ws.Range("D17").Value = "3.021"
$ws.Range("E17").Value = "-0.08%"

# Row 18 - BTSEToken
$ws.Range("D18").Value = "2.190"
$ws.Range("E18").Value = "8.32%"

# Row 19 - BitpandaEcosystemToken
$ws.Range("E19").Value = "-2.12%"

# Row 20 - LiechtensteinCryptoassetsExchange
$ws.Range("D20").Value = "0.03062"
$ws.Range("E20").Value = "-4.83%"

# Row 21 - ProBitToken
$ws.Range("E21").Value = "-1.50%"

# Row 22 - MCDex
$ws.Range("D22").Value = "3.759"
$ws.Range("E22").Value = "0.42%"

# Row 23 - CoinExToken
$ws.Range("D23").Value = "0.04580"
$ws.Range("E23").Value = "-2.48%"

# Row 25 - BitKan
$ws.Range("D25").Value = "0.001224"
$ws.Range("E25").Value = "-1.85%"

# Row 26 - HotbitToken
$ws.Range("D26").Value = "0.004492"
$ws.Range("E26").Value = "-2.52%"

# Row 27 - NitroEx
$ws.Range("D27").Value = "0.00009598"
$ws.Range("E27").Value = "-1.00%"

# Row 40 - IDEX
$ws.Range("D40").Value = "0.03639"
$ws.Range("E40").Value = "-0.47%"

# Row 41 - was KickToken, now BKEXToken
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "0.1371"
$ws.Range("E41").Value = "1.30%"

# Row 42 - was BKEXToken, now CEJI
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "0.002589"
$ws.Range("E42").Value = "2.41%"

# Row 43 - was CEJI, now KickToken
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "0.003451"
$ws.Range("E43").Value = "-43.95%"

# Row 44 - LocalTraders
$ws.Range("D44").Value = "0.008121"
$ws.Range("E44").Value = "4.48%"

# Row 45 - CoinLion
$ws.Range("D45").Value = "0.00005348"
$ws.Range("E45").Value = "0.64%"

# Row 46 - Kangarootoken
$ws.Range("E46").Value = "0.01%"

# Row 47 - CoinbaseStockToken
$ws.Range("D47").Value = "0.1090"
$ws.Range("E47").Value = "-18.34%"

# Row 48 - BOLO
$ws.Range("D48").Value = "0.002512"
$ws.Range("E48").Value = "18.23%"

# Row 49 - CryptobidCoin
$ws.Range("D49").Value = "0.00002100"
$ws.Range("E49").Value = "0.01%"

# Row 50 - SpecialPowerGold
$ws.Range("D50").Value = "0.0002000"
$ws.Range("E50").Value = "0.01%"

# Restore the default "Normal" style on the touched range so the
# temporary text number format doesn't leave a stray cell style behind.
$textRange.Style = "Normal"
